$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-10-03 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-04 Friday", 2) | Out-Null

# Update each arithmetic-problem cell in the table by (row, column) to avoid
# ambiguity from duplicate problem text (e.g. "92+1=93" appears twice).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "90-8=82"
$t.Cell(1, 2).Range.Text = "49+26=75"
$t.Cell(1, 3).Range.Text = "73+20=93"
$t.Cell(1, 4).Range.Text = "92-86=6"
$t.Cell(1, 5).Range.Text = "93-50=43"
$t.Cell(2, 1).Range.Text = "20+63=83"
$t.Cell(2, 2).Range.Text = "48-40=8"
$t.Cell(2, 3).Range.Text = "48+48=96"
$t.Cell(2, 4).Range.Text = "38+5=43"
$t.Cell(2, 5).Range.Text = "39-7=32"
$t.Cell(3, 1).Range.Text = "56-4=52"
$t.Cell(3, 2).Range.Text = "65-46=19"
$t.Cell(3, 3).Range.Text = "42-37=5"
$t.Cell(3, 4).Range.Text = "80-30=50"
$t.Cell(3, 5).Range.Text = "4+68=72"
$t.Cell(4, 1).Range.Text = "38-5=33"
$t.Cell(4, 2).Range.Text = "35+59=94"
$t.Cell(4, 3).Range.Text = "12+45=57"
$t.Cell(4, 4).Range.Text = "26+45=71"
$t.Cell(4, 5).Range.Text = "64+4=68"
$t.Cell(5, 1).Range.Text = "20+63=83"
$t.Cell(5, 2).Range.Text = "38+2=40"
$t.Cell(5, 3).Range.Text = "36+32=68"
$t.Cell(5, 4).Range.Text = "53+3=56"
$t.Cell(5, 5).Range.Text = "68-54=14"
$t.Cell(6, 1).Range.Text = "18+14=32"
$t.Cell(6, 2).Range.Text = "98-39=59"
$t.Cell(6, 3).Range.Text = "40-29=11"
$t.Cell(6, 4).Range.Text = "7+64=71"
$t.Cell(6, 5).Range.Text = "28+70=98"
$t.Cell(7, 1).Range.Text = "0+48=48"
$t.Cell(7, 2).Range.Text = "50+6=56"
$t.Cell(7, 3).Range.Text = "3+16=19"
$t.Cell(7, 4).Range.Text = "28+71=99"
$t.Cell(7, 5).Range.Text = "96-89=7"
$t.Cell(8, 1).Range.Text = "87-47=40"
$t.Cell(8, 2).Range.Text = "7+48=55"
$t.Cell(8, 3).Range.Text = "29+5=34"
$t.Cell(8, 4).Range.Text = "9+86=95"
$t.Cell(8, 5).Range.Text = "70-41=29"
$t.Cell(9, 1).Range.Text = "95-48=47"
$t.Cell(9, 2).Range.Text = "91-51=40"
$t.Cell(9, 3).Range.Text = "77-47=30"
$t.Cell(9, 4).Range.Text = "86+3=89"
$t.Cell(9, 5).Range.Text = "85-79=6"
$t.Cell(10, 1).Range.Text = "72-12=60"
$t.Cell(10, 2).Range.Text = "98-74=24"
$t.Cell(10, 3).Range.Text = "57-21=36"
$t.Cell(10, 4).Range.Text = "48-40=8"
$t.Cell(10, 5).Range.Text = "96-16=80"
$t.Cell(11, 1).Range.Text = "71-38=33"
$t.Cell(11, 2).Range.Text = "41+34=75"
$t.Cell(11, 3).Range.Text = "43-40=3"
$t.Cell(11, 4).Range.Text = "29+36=65"
$t.Cell(11, 5).Range.Text = "53-47=6"
$t.Cell(12, 1).Range.Text = "75-7=68"
$t.Cell(12, 2).Range.Text = "39-15=24"
$t.Cell(12, 3).Range.Text = "79-2=77"
$t.Cell(12, 4).Range.Text = "73+11=84"
$t.Cell(12, 5).Range.Text = "37+9=46"
$t.Cell(13, 1).Range.Text = "76-41=35"
$t.Cell(13, 2).Range.Text = "27+24=51"
$t.Cell(13, 3).Range.Text = "3+11=14"
$t.Cell(13, 4).Range.Text = "52-2=50"
$t.Cell(13, 5).Range.Text = "30+1=31"
$t.Cell(14, 1).Range.Text = "35+16=51"
$t.Cell(14, 2).Range.Text = "41-29=12"
$t.Cell(14, 3).Range.Text = "88-50=38"
$t.Cell(14, 4).Range.Text = "82-12=70"
$t.Cell(14, 5).Range.Text = "77-22=55"
$t.Cell(15, 1).Range.Text = "60-25=35"
$t.Cell(15, 2).Range.Text = "78-71=7"
$t.Cell(15, 3).Range.Text = "33+36=69"
$t.Cell(15, 4).Range.Text = "90-28=62"
$t.Cell(15, 5).Range.Text = "44+24=68"
$t.Cell(16, 1).Range.Text = "91-26=65"
$t.Cell(16, 2).Range.Text = "4+5=9"
$t.Cell(16, 3).Range.Text = "79-44=35"
$t.Cell(16, 4).Range.Text = "53+32=85"
$t.Cell(16, 5).Range.Text = "53+9=62"
$t.Cell(17, 1).Range.Text = "85-16=69"
$t.Cell(17, 2).Range.Text = "33+55=88"
$t.Cell(17, 3).Range.Text = "27+4=31"
$t.Cell(17, 4).Range.Text = "93-9=84"
$t.Cell(17, 5).Range.Text = "93-72=21"
$t.Cell(18, 1).Range.Text = "60+4=64"
$t.Cell(18, 2).Range.Text = "25+23=48"
$t.Cell(18, 3).Range.Text = "7+1=8"
$t.Cell(18, 4).Range.Text = "73-24=49"
$t.Cell(18, 5).Range.Text = "13+80=93"
$t.Cell(19, 1).Range.Text = "79-41=38"
$t.Cell(19, 2).Range.Text = "62-27=35"
$t.Cell(19, 3).Range.Text = "48+27=75"
$t.Cell(19, 4).Range.Text = "80-17=63"
$t.Cell(19, 5).Range.Text = "76-51=25"
$t.Cell(20, 1).Range.Text = "71+22=93"
$t.Cell(20, 2).Range.Text = "53+32=85"
$t.Cell(20, 3).Range.Text = "89-58=31"
$t.Cell(20, 4).Range.Text = "48-5=43"
$t.Cell(20, 5).Range.Text = "38+14=52"

Write-Output "done"
